# Generate Report for Handoff
#
# Four files (0f81beb4, a2d9c3fc, b91a9bab, ba3ce974) just finished their
# handoff xliff generation, so:
#   - Overview!G (Latest HO Xliff Generate Date) advances from 10:35:33 to 10:35:48
#     for those four rows (4-7)
#   - zh-cn / de-de Priority (col E) flips from "low" to "ht" for rows 4-7
#   - zh-cn / de-de Latest Handoff Datetime (col H) advances for rows 4-7
#     (zh-cn: 10:35:28 -> 10:35:43 ; de-de: 10:35:33 -> 10:35:48)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

foreach ($row in 4..7) {
    $wsOverview.Range("G$row").Value = "2016-08-18 10:35:48"

    $wsZhCn.Range("E$row").Value = "ht"
    $wsZhCn.Range("H$row").Value = "2016-08-18 10:35:43"

    $wsDeDe.Range("E$row").Value = "ht"
    $wsDeDe.Range("H$row").Value = "2016-08-18 10:35:48"
}
